$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# New header cell N1 = "Correction", matching style of existing header M1
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "Correction"

# Fill previously-empty M column cells with "nan" (rows 2-12)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# Add new empty-string cells in column N for rows 2-12
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = ""
}
